$d = $word.ActiveDocument
$find = $d.Content.Find

# --- Merge runs that were split by spelling/grammar proofing marks back
# --- into single plain runs (no visible text change, just de-fragmenting
# --- the <w:r>/<w:proofErr> markup Word's proofer had introduced).

$p1 = "We hold these truths to be self-evident; that all men and women are created equal; that they are endowed by their Creator with certain inalienable rights; that among these are life, liberty, and the pursuit of happiness; that to secure these rights governments are instituted, deriving their just powers from the consent of the governed. Whenever any form of government becomes destructive of these ends, it is the right of those who suffer from it to refuse allegiance to it, and to insist upon the institution of a new government, laying its foundation on such principles, and organizing its powers in such form, as to them shall seem most likely to effect their safety and happiness. Prudence, indeed, will dictate that governments long established should not be changed for light and transient causes; and, accordingly, all experience hath shown that mankind are more disposed to suffer, while evils are sufferable, than to right themselves by abolishing the forms to which they were accustomed. But when a long train of abuses and usurpations, pursuing invariably the same object, evinces a design to reduce them under absolute despotism, it is their duty to throw off such government, and to provide new guards for their future security. Such has been the patient sufferance of the women under this government, and such is now the necessity which constrains them to demand the equal station to which they are entitled."
$find.Execute($p1, $true, $false, $false, $false, $false, $true, 1, $false, $p1, 2) | Out-Null

$p2 = "He has made her morally, an irresponsible being, as she can commit many crimes with impunity, provided they be done in the presence of her husband. In the covenant of marriage, she is compelled to promise obedience to her husband, he becoming, to all intents and purposes, her master - the law giving him power to deprive her of her liberty, and to administer chastisement."
$find.Execute($p2, $true, $false, $false, $false, $false, $true, 1, $false, $p2, 2) | Out-Null

$p3 = "In entering upon the great work before us, we anticipate no small amount of misconception, misrepresentation, and ridicule; but we shall use every instrumentality within our power to effect our object. We shall employ agents, circulate tracts, petition the State and national Legislatures, and endeavor to enlist the pulpit and the press in our behalf. We hope this Convention will be followed by a series of Conventions, embracing every part of the country."
$find.Execute($p3, $true, $false, $false, $false, $false, $true, 1, $false, $p3, 2) | Out-Null

# --- Add a second trailing empty paragraph at the end of the document
# --- (there was one empty paragraph before the sectPr; now there are two).

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

# --- Give the "Normal" style explicit paragraph spacing
# --- (w:spacing w:line="256" w:lineRule="auto") instead of inheriting
# --- the implicit Word default.

$normalStyle = $d.Styles("Normal")
$normalStyle.ParagraphFormat.LineSpacingRule = 5
$normalStyle.ParagraphFormat.LineSpacing = 12.8
